# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows that changed after re-pulling data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 1
    4  = -1
    5  = -3
    7  = -3
    8  = 1
    9  = -5
    11 = -6
    12 = 1
    13 = 2
    14 = 3
    15 = -1
    16 = 5
    20 = 7
    21 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
